# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is refreshed for zh-cn and de-de
#  - The stale "handback file is not latest" Error Detail message is cleared

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: per-locale status columns (zh-cn, de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-26 15:02:03"
$wsZhCn.Range("P2").Value = ""

# de-de detail sheet
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-26 15:02:20"
$wsDeDe.Range("P2").Value = ""

# Column widths auto-adjust after the report regeneration: the Status columns
# grow to fit the longer "Handed back..." text, and the now-empty Error Detail
# columns shrink back down to fit their (now blank) content.
$statusColWidth = 29.166666666666664   # widened Status column (was ~17.2)
$errorColWidth  = 12.833333333333332   # narrowed Error Detail column (was 40)

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth   # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth       # C: Status
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColWidth       # P: Error Detail
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth       # C: Status
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColWidth       # P: Error Detail
